$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row heights (rows 1-3 get new explicit/custom heights) ---
$ws.Rows.Item(1).RowHeight = 66
$ws.Rows.Item(2).RowHeight = 14.25
$ws.Rows.Item(3).RowHeight = 14.25

# --- New column N (14) data, row by row, matching the format of column M in the same row ---

# Row 3: empty cell, bottom-medium-border style (same as A3:M3)
$ws.Range("N3").Borders.Item(9).LineStyle = -4138
$ws.Range("N3").Borders.Item(9).Weight = -4138
$ws.Range("N3").Borders.Item(9).Color = 0

# Row 4: year header 2023, bold Times New Roman 9, right aligned, vertical centered, bottom border
$c = $ws.Range("N4")
$c.Value = 2023
$c.Font.Name = "Times New Roman"
$c.Font.Size = 9
$c.Font.Bold = $true
$c.HorizontalAlignment = -4152
$c.VerticalAlignment = -4108
$c.Borders.Item(9).LineStyle = -4138
$c.Borders.Item(9).Weight = -4138
$c.Borders.Item(9).Color = 0

# Row 5: value 0, numeric "0.0" format, right aligned, vertical centered
$c = $ws.Range("N5")
$c.Value = 0
$c.Font.Name = "Times New Roman"
$c.Font.Size = 9
$c.NumberFormat = "0.0"
$c.HorizontalAlignment = -4152
$c.VerticalAlignment = -4108

# Row 6: value 48.5, numeric "0.0" format, right aligned (no vertical center; font without theme color)
$c = $ws.Range("N6")
$c.Value = 48.5
$c.Font.Name = "Times New Roman"
$c.Font.Size = 9
$c.NumberFormat = "0.0"
$c.HorizontalAlignment = -4152

# Row 7: value 23.2, General number format, vertical centered
$c = $ws.Range("N7")
$c.Value = 23.2
$c.Font.Name = "Times New Roman"
$c.Font.Size = 9
$c.NumberFormat = "General"
$c.VerticalAlignment = -4108

# Row 8: value 19.3, numeric "0.0" format, vertical centered
$c = $ws.Range("N8")
$c.Value = 19.3
$c.Font.Name = "Times New Roman"
$c.Font.Size = 9
$c.NumberFormat = "0.0"
$c.VerticalAlignment = -4108

# Row 9: value 9.1, numeric "0.0" format, right aligned, vertical centered, bottom border
$c = $ws.Range("N9")
$c.Value = 9.1
$c.Font.Name = "Times New Roman"
$c.Font.Size = 9
$c.NumberFormat = "0.0"
$c.HorizontalAlignment = -4152
$c.VerticalAlignment = -4108
$c.Borders.Item(9).LineStyle = -4138
$c.Borders.Item(9).Weight = -4138
$c.Borders.Item(9).Color = 0

# --- L7: change display format to one-decimal ("0.0") to match new style ---
$ws.Range("L7").NumberFormat = "0.0"

# --- Reset the saved selection back to the default top-left cell ---
$ws.Range("A1").Select()
